$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 29, shifting existing rows 29:40 down to 30:41.
$ws.Rows("29:29").Insert()

# Populate the new row 29 with the new weekly record.
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C29").Value = "Arica y Parinacota"
$ws.Range("D29").Value = 44762
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
$ws.Range("E29").Value = 15
$ws.Range("F29").Value = 100112052
$ws.Range("G29").Value = "Albahaca"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 1800
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1900
$ws.Range("N29").Value = "$/paquete"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 1900
$ws.Range("Q29").Value = 1
$ws.Range("R29").Value = "Hortaliza"
